$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.598.76'
$ws.Range("E2").Value = '  +1.14%  '

$ws.Range("D3").Value = '3.318.82'
$ws.Range("E3").Value = '  +1.95%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  +0.15%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '580.55'
$ws.Range("E5").Value = '  +0.35%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '174.08'
$ws.Range("E6").Value = '  +0.53%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.587'
$ws.Range("E8").Value = '  +2.11%  '

$ws.Range("D9").Value = '3.315.03'

$ws.Range("E10").Value = '  +5.29%  '

$ws.Range("E11").Value = '  +1.62%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '46.63'
$ws.Range("E12").Value = '  +4.30%  '

$ws.Range("E13").Value = '  +0.62%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '692.38'
$ws.Range("E14").Value = '  +3.97%  '

$ws.Range("D15").Value = '3.857.89'
$ws.Range("E15").Value = '  +2.19%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '8.35'
$ws.Range("E16").Value = '  +1.53%  '

$ws.Range("D17").Value = '67.556.03'
$ws.Range("E17").Value = '  +1.26%  '

$ws.Range("E18").Value = '  +0.67%  '

$ws.Range("D19").Value = '3.316.71'
$ws.Range("E19").Value = '  +2.23%  '

$ws.Range("E20").Value = '  +1.55%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '11.07'
$ws.Range("E21").Value = '  +3.29%  '

$ws.Range("E22").Value = '  +1.51%  '

$ws.Range("E23").Value = '  +2.92%  '

$ws.Range("E24").Value = '  -0.50%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '101.42'
$ws.Range("E25").Value = '  +5.06%  '

$ws.Range("E26").Value = '  +1.39%  '

$ws.Range("E27").Value = '  +1.59%  '

$ws.Range("E28").Value = '  +2.96%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '32.79'
$ws.Range("E29").Value = '  +2.43%  '

$ws.Range("E30").Value = '  +2.41%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.97'
$ws.Range("E31").Value = '  +2.23%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '567.76'
$ws.Range("E32").Value = '  +0.75%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '10.96'
$ws.Range("E33").Value = '  +0.90%  '

$ws.Range("E34").Value = '  +3.00%  '

$ws.Range("E35").Value = '  +0.00%  '

$ws.Range("B36").Value = 'OKB'
$ws.Range("C36").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '57.19'
$ws.Range("E36").Value = '  +3.16%  '

$ws.Range("B37").Value = 'Maker'
$ws.Range("C37").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D37").Value = '3.708.23'
$ws.Range("E37").Value = '  -1.10%  '

$ws.Range("E38").Value = '  -5.79%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '34.90'
$ws.Range("E39").Value = '  +8.26%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.132'
$ws.Range("E40").Value = '  +2.71%  '

$ws.Range("E41").Value = '  +4.68%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.60'
$ws.Range("E42").Value = '  +0.01%  '

$ws.Range("B43").Value = 'ApeXProtocol'
$ws.Range("C43").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.31'
$ws.Range("E43").Value = '  +3.03%  '

$ws.Range("B44").Value = 'TheGraph'
$ws.Range("C44").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.333'
$ws.Range("E44").Value = '  +2.77%  '

$ws.Range("E45").Value = '  +1.18%  '

$ws.Range("E46").Value = '  +2.19%  '

$ws.Range("E47").Value = '  +3.45%  '

$ws.Range("E48").Value = '  +1.47%  '

$ws.Range("E49").Value = '  -0.01%  '

$ws.Range("E50").Value = '  -1.02%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '131.21'
$ws.Range("E51").Value = '  +2.16%  '
